# BIS-769: Fixed xls test files
# Adds "Pattern" / "Pattern Type" columns (M/N) to the two property-definition
# header tables (rows 12 and 20) of the sample-type export/import test sheet,
# copying the formatting of the existing "Unique" header cell (column L) so the
# new header cells render identically, then moves the active selection to the
# newly added range, matching the upstream fixture update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-PatternColumns($headerRow) {
    $lCell = $ws.Range("L" + $headerRow)
    $mCell = $ws.Range("M" + $headerRow)
    $nCell = $ws.Range("N" + $headerRow)

    # Clone the "Unique" header cell's formatting onto the two new cells.
    $lCell.Copy()
    $mCell.PasteSpecial(-4122)
    $nCell.PasteSpecial(-4122)

    $mCell.Value = "Pattern"
    $nCell.Value = "Pattern Type"
}

# First property table (ANTIBODY sample type), header row 12.
Add-PatternColumns 12

# Second property table (VIRUS sample type), header row 20.
Add-PatternColumns 20

# Move the selection to the newly added cells, as in the updated fixture.
$ws.Range("M20:N20").Select() | Out-Null
